$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Clear the header cell H1 (was "DATA_REPERFORMING") -- the column now represents
# a reperforming flag rather than a reperforming date, so the header text is removed
# (style/formatting on the cell is preserved).
$ws.Range("H1").ClearContents()

# Move the active cell selection from H2 to H1
$ws.Range("H1").Select()
